$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-03-04 Tuesday" "2025-03-05 Wednesday"

Replace-Text "943×9=8487" "567×8=4536"
Replace-Text "828×8=6624" "882×6=5292"
Replace-Text "684×4=2736" "793×2=1586"
Replace-Text "195×3=585" "643×8=5144"
Replace-Text "212×4=848" "449×2=898"
Replace-Text "522×6=3132" "146×5=730"
Replace-Text "362×6=2172" "961×9=8649"
Replace-Text "535×9=4815" "367×7=2569"
Replace-Text "590×8=4720" "645×5=3225"
Replace-Text "527×5=2635" "504×6=3024"
Replace-Text "477×9=4293" "433×8=3464"
Replace-Text "213×9=1917" "355×2=710"
Replace-Text "237×9=2133" "857×4=3428"
Replace-Text "711×6=4266" "295×6=1770"
Replace-Text "899×3=2697" "634×3=1902"
Replace-Text "227×2=454" "725×3=2175"
Replace-Text "352×4=1408" "860×5=4300"
Replace-Text "321×9=2889" "689×2=1378"
Replace-Text "953×3=2859" "384×4=1536"
Replace-Text "964×9=8676" "499×3=1497"
Replace-Text "240×5=1200" "330×3=990"
Replace-Text "287×3=861" "400×5=2000"
Replace-Text "843×6=5058" "101×6=606"
Replace-Text "696×3=2088" "657×6=3942"
Replace-Text "699×7=4893" "113×2=226"
